# Fix typo in "Lettre Recommandée avec Accusé de Réception" -> lowercase the
# adjective/participle parts: "Lettre recommandée avec accusé de réception".
# Each target word begins a separate run in the source document, so we match
# and replace them individually to keep the edits narrowly scoped.

$d = $word.ActiveDocument

$d.Content.Find.Execute("Lettre Recommand", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Lettre recommand", 2) | Out-Null

$d.Content.Find.Execute("e avec Accus", $true, $false, $false, $false, $false, `
    $true, 1, $false, "e avec accus", 2) | Out-Null

$d.Content.Find.Execute("de R", $true, $false, $false, $false, $false, `
    $true, 1, $false, "de r", 2) | Out-Null

# Replace the straight apostrophe with a typographic (curly) one in
# "l'expression de ma considération".
$d.Content.Find.Execute("Monsieur, l'expression", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Monsieur, l’expression", 2) | Out-Null
